{"js": "// Apply Papiamento wording fixes to the \"Children with Disabilities\" facilitator\n// document. Each entry below is an exact, whole-run-text replacement taken\n// from the canonical OOXML diff: we locate the unique original run text in\n// the document body and swap it for the corrected text, leaving every other\n// run (including bold labels like \"Aktitutnan:\") and all paragraph/run\n// formatting completely untouched.\nconst replacements = [\n  {\n    \"old\": \"Mi\u00e9ntras ku ta konsidera abla e forma prinsipal di komunikashon, palabra no ta e \u00faniko manera pa transmit\u00ed nifikashon. Lenguahe di kurpa, ekspreshon di kara i tono di bos tambe ta medionan poderoso di komunikashon. Esaki tin implikashon pa muchanan ku diferente desabilidat.    \",\n    \"new\": \"Mi\u00e9ntras ku ta konsider\u00e1 Papia e forma prinsipal di komunikashon, palabra no ta e \u00faniko manera pa transmit\u00ed nifikashon. Lenguahe di kurpa, ekspreshon di kara i tono di bos tambe ta medionan poderoso di komunikashon. Esaki tin implikashon pa muchanan ku diferente desabilidat.    \"\n  },\n  {\n    \"old\": \"Tur hende, sin import\u00e1 nan desabilidat, por komunik\u00e1 di algun manera. Sinembargo, muchanan ku desabilidat por enfrent\u00e1 bareranan signifikante pa komunik\u00e1 nan punto di bista i sintimentunan. Un di e bareranan prinsipal ta ku hendenan sin desabilidat hopi bia no tin un kompromiso pa komunik\u00e1 ku mucha \u00f2f adultonan ku desabilidat. Pero ta posibel pa super\u00e1 e bareranan ak\u00ed.  \",\n    \"new\": \"Tur hende, sin import\u00e1 nan desabilidat, por komunik\u00e1 di algun manera. Sinembargo, muchanan ku desabilidat por enfrent\u00e1 bareranan signifikante pa komunik\u00e1 nan punto di bista i sintimentunan. Un di e bareranan prinsipal ta ku hendenan sin desabilidat hopi bia no ta kompremet\u00e9 nan mes pa komunik\u00e1 ku mucha \u00f2f adultonan ku desabilidat. Pero ta posibel pa super\u00e1 e bareranan ak\u00ed.  \"\n  },\n  {\n    \"old\": \" Hopi ta kere ku muchanan ku limitashon ku no ta komunika di e mesun maneranan ku otro muchanan no ta inteligente \u00f2f ta desobediente. E aktitut negativo aki ta un barera grandi pa komunikashon efektivo. \",\n    \"new\": \" Hopi ta kere ku muchanan ku desabilidat ku no ta komunik\u00e1 di e mesun manera ku otro muchanan no ta inteligente \u00f2f ta desobediente. E aktitut negativo aki ta un barera grandi pa komunikashon efektivo. \"\n  },\n  {\n    \"old\": \" Normalmente ta mas dif\u00edsil pa mucha muh\u00e9 ku desabilidat ekspres\u00e1 nan punto di bista. Den mayoria kultura, e ekspektativa muchu mas grandi di pasividat di mucha muh\u00e9, ta krea e bareranan ku nan ta enfrent\u00e1 deb\u00ed na desabilidatnan. \",\n    \"new\": \" Normalmente ta mas dif\u00edsil pa mucha muh\u00e9 ku desabilidat ekspres\u00e1 nan punto di bista. Den mayoria kultura, e ekspektativa grandi di pasividat di mucha muh\u00e9, ta krea e bareranan ku nan ta enfrent\u00e1 deb\u00ed na desabilidatnan. \"\n  },\n  {\n    \"old\": \" Muchanan ku desabilidat ku no ta kustum\u00e1 ku ta puntra nan nan opinion \u00f2f ku no ta kustuma ku hende ta skucha nan, mester di mas tempu pa krea konfiansa i seguridat. E lo por tuma tempu pa eksplor\u00e1 e mih\u00f3 maneranan di komunik\u00e1 ku un mucha en partikular. Mayornan lo mester tin mas pasenshi ku nan mes i nan yunan ora di pasa Tempu Huntu ku nan yunan. \",\n    \"new\": \" Muchanan ku desabilidat ku no ta kustumbr\u00e1 ku ta puntra nan nan opinion \u00f2f ku no ta kustumbr\u00e1 ku hende ta skucha nan, mester di mas tempu pa krea konfiansa i seguridat. E lo por tuma tempu pa eksplor\u00e1 e mih\u00f3 maneranan di komunik\u00e1 ku un mucha en partikular. Mayornan lo mester tin mas pasenshi ku nan mes i nan yunan ora di pasa Tempu Huntu ku nan yunan. \"\n  },\n  {\n    \"old\": \"Ta tuma tempu pa muchanan ku desabilidat por krea konfiansa i seguridat, ya ku nan no ta kustum\u00e1 ku ta puntra nan na opinion \u00f2f ku ta skucha nan.\",\n    \"new\": \"Ta tuma tempu pa muchanan ku desabilidat por krea konfiansa den nan mes i konfiansa den otro, ya ku nan no ta kustumbr\u00e1 ku ta puntra nan na opinion \u00f2f ku ta skucha nan.\"\n  },\n  {\n    \"old\": \"Si kualke hende den e famia tin desabilidat, s\u00f2ru pa nan ta igualmente enbolb\u00ed. Mucha i adultonan ku desabilidat por enfrent\u00e1 bareranan signifikante pa komunik\u00e1 nan puntonan di bista i sintimentunan. Un di e bareranan prinsipal ta ku hendenan sin desabilidat hopi bia no tin un kompromiso pa komunik\u00e1 ku mucha \u00f2f adultonan ku desabilidat. Pero ta posibel pa super\u00e1 e bareranan ak\u00ed. Tin bia e proseso ak\u00ed por tuma hopi tempu, pero ta nesesario pa garantis\u00e1 nan derecho di partisip\u00e1. \",\n    \"new\": \"Si kualke hende den e famia tin desabilidat, s\u00f2ru pa nan ta igualmente enbolb\u00ed. Mucha i adultonan ku desabilidat por enfrent\u00e1 bareranan signifikante pa komunik\u00e1 nan puntonan di bista i sintimentunan. Un di e bareranan prinsipal ta ku hendenan sin desabilidat hopi bia no ta kompromet\u00e9 nan mes pa komunik\u00e1 ku mucha \u00f2f adultonan ku desabilidat. Pero ta posibel pa super\u00e1 e bareranan ak\u00ed. Tin bia e proseso ak\u00ed por tuma hopi tempu, pero ta nesesario pa garantis\u00e1 nan derecho di partisip\u00e1. \"\n  }\n];\n\nfor (const { old, new: replacement } of replacements) {\n  const results = context.document.body.search(old, { matchCase: true, matchWholeWord: false });\n  results.load(\"items/text\");\n  await context.sync();\n\n  if (results.items.length !== 1) {\n    throw new Error(\n      `Expected exactly 1 match for \"${old.slice(0, 40)}...\" but found ${results.items.length}`\n    );\n  }\n\n  results.items[0].insertText(replacement, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Apply Papiamento wording fixes to the \"Children with Disabilities\" facilitator\n# document. Each entry below is an exact, whole-run-text replacement taken\n# from the canonical OOXML diff: Find/Replace locates the unique original\n# run text in the document and swaps it for the corrected text, leaving\n# every other run (including bold labels like \"Aktitutnan:\") and all\n# paragraph/run formatting completely untouched.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"Mi\u00e9ntras ku ta konsidera abla e forma prinsipal di komunikashon, palabra no ta e \u00faniko manera pa transmit\u00ed nifikashon. Lenguahe di kurpa, ekspreshon di kara i tono di bos tambe ta medionan poderoso di komunikashon. Esaki tin implikashon pa muchanan ku diferente desabilidat.    \"; New = \"Mi\u00e9ntras ku ta konsider\u00e1 Papia e forma prinsipal di komunikashon, palabra no ta e \u00faniko manera pa transmit\u00ed nifikashon. Lenguahe di kurpa, ekspreshon di kara i tono di bos tambe ta medionan poderoso di komunikashon. Esaki tin implikashon pa muchanan ku diferente desabilidat.    \" },\n    @{ Old = \"Tur hende, sin import\u00e1 nan desabilidat, por komunik\u00e1 di algun manera. Sinembargo, muchanan ku desabilidat por enfrent\u00e1 bareranan signifikante pa komunik\u00e1 nan punto di bista i sintimentunan. Un di e bareranan prinsipal ta ku hendenan sin desabilidat hopi bia no tin un kompromiso pa komunik\u00e1 ku mucha \u00f2f adultonan ku desabilidat. Pero ta posibel pa super\u00e1 e bareranan ak\u00ed.  \"; New = \"Tur hende, sin import\u00e1 nan desabilidat, por komunik\u00e1 di algun manera. Sinembargo, muchanan ku desabilidat por enfrent\u00e1 bareranan signifikante pa komunik\u00e1 nan punto di bista i sintimentunan. Un di e bareranan prinsipal ta ku hendenan sin desabilidat hopi bia no ta kompremet\u00e9 nan mes pa komunik\u00e1 ku mucha \u00f2f adultonan ku desabilidat. Pero ta posibel pa super\u00e1 e bareranan ak\u00ed.  \" },\n    @{ Old = \" Hopi ta kere ku muchanan ku limitashon ku no ta komunika di e mesun maneranan ku otro muchanan no ta inteligente \u00f2f ta desobediente. E aktitut negativo aki ta un barera grandi pa komunikashon efektivo. \"; New = \" Hopi ta kere ku muchanan ku desabilidat ku no ta komunik\u00e1 di e mesun manera ku otro muchanan no ta inteligente \u00f2f ta desobediente. E aktitut negativo aki ta un barera grandi pa komunikashon efektivo. \" },\n    @{ Old = \" Normalmente ta mas dif\u00edsil pa mucha muh\u00e9 ku desabilidat ekspres\u00e1 nan punto di bista. Den mayoria kultura, e ekspektativa muchu mas grandi di pasividat di mucha muh\u00e9, ta krea e bareranan ku nan ta enfrent\u00e1 deb\u00ed na desabilidatnan. \"; New = \" Normalmente ta mas dif\u00edsil pa mucha muh\u00e9 ku desabilidat ekspres\u00e1 nan punto di bista. Den mayoria kultura, e ekspektativa grandi di pasividat di mucha muh\u00e9, ta krea e bareranan ku nan ta enfrent\u00e1 deb\u00ed na desabilidatnan. \" },\n    @{ Old = \" Muchanan ku desabilidat ku no ta kustum\u00e1 ku ta puntra nan nan opinion \u00f2f ku no ta kustuma ku hende ta skucha nan, mester di mas tempu pa krea konfiansa i seguridat. E lo por tuma tempu pa eksplor\u00e1 e mih\u00f3 maneranan di komunik\u00e1 ku un mucha en partikular. Mayornan lo mester tin mas pasenshi ku nan mes i nan yunan ora di pasa Tempu Huntu ku nan yunan. \"; New = \" Muchanan ku desabilidat ku no ta kustumbr\u00e1 ku ta puntra nan nan opinion \u00f2f ku no ta kustumbr\u00e1 ku hende ta skucha nan, mester di mas tempu pa krea konfiansa i seguridat. E lo por tuma tempu pa eksplor\u00e1 e mih\u00f3 maneranan di komunik\u00e1 ku un mucha en partikular. Mayornan lo mester tin mas pasenshi ku nan mes i nan yunan ora di pasa Tempu Huntu ku nan yunan. \" },\n    @{ Old = \"Ta tuma tempu pa muchanan ku desabilidat por krea konfiansa i seguridat, ya ku nan no ta kustum\u00e1 ku ta puntra nan na opinion \u00f2f ku ta skucha nan.\"; New = \"Ta tuma tempu pa muchanan ku desabilidat por krea konfiansa den nan mes i konfiansa den otro, ya ku nan no ta kustumbr\u00e1 ku ta puntra nan na opinion \u00f2f ku ta skucha nan.\" },\n    @{ Old = \"Si kualke hende den e famia tin desabilidat, s\u00f2ru pa nan ta igualmente enbolb\u00ed. Mucha i adultonan ku desabilidat por enfrent\u00e1 bareranan signifikante pa komunik\u00e1 nan puntonan di bista i sintimentunan. Un di e bareranan prinsipal ta ku hendenan sin desabilidat hopi bia no tin un kompromiso pa komunik\u00e1 ku mucha \u00f2f adultonan ku desabilidat. Pero ta posibel pa super\u00e1 e bareranan ak\u00ed. Tin bia e proseso ak\u00ed por tuma hopi tempu, pero ta nesesario pa garantis\u00e1 nan derecho di partisip\u00e1. \"; New = \"Si kualke hende den e famia tin desabilidat, s\u00f2ru pa nan ta igualmente enbolb\u00ed. Mucha i adultonan ku desabilidat por enfrent\u00e1 bareranan signifikante pa komunik\u00e1 nan puntonan di bista i sintimentunan. Un di e bareranan prinsipal ta ku hendenan sin desabilidat hopi bia no ta kompromet\u00e9 nan mes pa komunik\u00e1 ku mucha \u00f2f adultonan ku desabilidat. Pero ta posibel pa super\u00e1 e bareranan ak\u00ed. Tin bia e proseso ak\u00ed por tuma hopi tempu, pero ta nesesario pa garantis\u00e1 nan derecho di partisip\u00e1. \" }\n)\n\nforeach ($item in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n\n    $found = $find.Execute($item.Old, $false, $false, $false, $false, $false, $true, 1, $false, $item.New, 2)\n\n    if (-not $found) {\n        throw \"Find/Replace did not find expected text: $($item.Old.Substring(0, 40))...\"\n    }\n}\n"}
